$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Refresh the per-row query timestamps on the "data" sheet (F2:F74) to
# reflect the re-run of the PanelApp query.
$timestamps = @(
    "2021-10-05 14:20:53.313660",
    "2021-10-05 14:20:53.313668",
    "2021-10-05 14:20:53.313671",
    "2021-10-05 14:20:53.313674",
    "2021-10-05 14:20:53.313676",
    "2021-10-05 14:20:53.313679",
    "2021-10-05 14:20:53.313681",
    "2021-10-05 14:20:53.313684",
    "2021-10-05 14:20:53.313686",
    "2021-10-05 14:20:53.313689",
    "2021-10-05 14:20:53.313691",
    "2021-10-05 14:20:53.313693",
    "2021-10-05 14:20:53.313696",
    "2021-10-05 14:20:53.313698",
    "2021-10-05 14:20:53.313701",
    "2021-10-05 14:20:53.313703",
    "2021-10-05 14:20:53.313706",
    "2021-10-05 14:20:53.313708",
    "2021-10-05 14:20:53.313711",
    "2021-10-05 14:20:53.313713",
    "2021-10-05 14:20:53.313715",
    "2021-10-05 14:20:53.313718",
    "2021-10-05 14:20:53.313720",
    "2021-10-05 14:20:53.313722",
    "2021-10-05 14:20:53.313725",
    "2021-10-05 14:20:53.313728",
    "2021-10-05 14:20:53.313730",
    "2021-10-05 14:20:53.313733",
    "2021-10-05 14:20:53.313735",
    "2021-10-05 14:20:53.313738",
    "2021-10-05 14:20:53.313740",
    "2021-10-05 14:20:53.313742",
    "2021-10-05 14:20:53.313745",
    "2021-10-05 14:20:53.313748",
    "2021-10-05 14:20:53.313750",
    "2021-10-05 14:20:53.313753",
    "2021-10-05 14:20:53.313755",
    "2021-10-05 14:20:53.313757",
    "2021-10-05 14:20:53.313760",
    "2021-10-05 14:20:53.313762",
    "2021-10-05 14:20:53.313765",
    "2021-10-05 14:20:53.313767",
    "2021-10-05 14:20:53.313769",
    "2021-10-05 14:20:53.313772",
    "2021-10-05 14:20:53.313774",
    "2021-10-05 14:20:53.313776",
    "2021-10-05 14:20:53.313779",
    "2021-10-05 14:20:53.313781",
    "2021-10-05 14:20:53.313783",
    "2021-10-05 14:20:53.313786",
    "2021-10-05 14:20:53.313788",
    "2021-10-05 14:20:53.313790",
    "2021-10-05 14:20:53.313793",
    "2021-10-05 14:20:53.313796",
    "2021-10-05 14:20:53.313798",
    "2021-10-05 14:20:53.313800",
    "2021-10-05 14:20:53.313803",
    "2021-10-05 14:20:53.313805",
    "2021-10-05 14:20:53.313808",
    "2021-10-05 14:20:53.313813",
    "2021-10-05 14:20:53.313816",
    "2021-10-05 14:20:53.313818",
    "2021-10-05 14:20:53.313821",
    "2021-10-05 14:20:53.313823",
    "2021-10-05 14:20:53.313827",
    "2021-10-05 14:20:53.313829",
    "2021-10-05 14:20:53.313832",
    "2021-10-05 14:20:53.313834",
    "2021-10-05 14:20:53.313836",
    "2021-10-05 14:20:53.313839",
    "2021-10-05 14:20:53.313841",
    "2021-10-05 14:20:53.313844",
    "2021-10-05 14:20:53.313846"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

# Add a new "metadata" tab (sheetId 2), positioned after "data".
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

# Header row.
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Match the bold/bordered/centered header formatting used on the "data" tab.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row describing the panel query that produced the "data" tab.
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Hypertrophic cardiomyopathy - teen and adult"
$meta.Range("C2").Value = 49

$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.26"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2021-09-28T09:21:35.238412Z"
$meta.Range("F2").Value = "2021-10-05 14:20:53.310731"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/49/?format=json"

# Match the row-index formatting used in column A on the "data" tab.
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
